# Update odds values on the active worksheet to match the latest
# FlashScore data refresh ("Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("V2").Value = 1.57

# Row 3
$ws.Range("J3").Value = 2.38
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8

# Row 4
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 2.88

# Row 5
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5

# Row 6
$ws.Range("G6").Value = 5.25
$ws.Range("N6").Value = 8
$ws.Range("R6").Value = 1.6
$ws.Range("X6").Value = 23
$ws.Range("Y6").Value = 17
$ws.Range("AE6").Value = 19
$ws.Range("AQ6").Value = 101
$ws.Range("AW6").Value = 3.6

# Row 7
$ws.Range("G7").Value = 1.8
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 2.5
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5
$ws.Range("R7").Value = 1.53
$ws.Range("X7").Value = 7
$ws.Range("Z7").Value = 13
$ws.Range("AH7").Value = 10
$ws.Range("AI7").Value = 23
$ws.Range("AQ7").Value = 34
$ws.Range("AT7").Value = 2.38
$ws.Range("AW7").Value = 6.5
$ws.Range("AZ7").Value = 126

# Row 9
$ws.Range("R9").Value = 1.48

# Row 12
$ws.Range("J12").Value = 2.25
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("O12").Value = 1.4
$ws.Range("P12").Value = 2.75
$ws.Range("U12").Value = 2.25
$ws.Range("V12").Value = 1.57
$ws.Range("W12").Value = 5.5
$ws.Range("AC12").Value = 7
$ws.Range("AE12").Value = 21
$ws.Range("AK12").Value = 67
$ws.Range("AP12").Value = 23

# Row 13
$ws.Range("Q13").Value = 1.98
$ws.Range("R13").Value = 1.88
